$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jurado 2 for rows 8 and 15, and Jurado 1 for row 21 previously referenced
# "NOBLECILLA VINCES WILLIAM ALFREDO" -- update them to the shortened name
# "NOBLECILLA VINCES WILLIAM" (adds a new shared string, reused by all three).
$ws.Range("G8").Value = "NOBLECILLA VINCES WILLIAM"
$ws.Range("G15").Value = "NOBLECILLA VINCES WILLIAM"

# F21 also gets a distinct font color (explicit black) applied along with the
# same text update, producing a new font + cell style.
$ws.Range("F21").Font.Color = 0
$ws.Range("F21").Value = "NOBLECILLA VINCES WILLIAM"

# Update the active selection to reflect where the edit left the cursor.
[void]$ws.Range("F21").Select()

# Page orientation explicitly set to portrait.
$ws.PageSetup.Orientation = 1
